# Apply the cryptos list update (prices/volumes refreshed, rows 25/26 re-ranked).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "57.705.70"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "2.444.79"
$ws.Range("E3").Value = "  -3.68%  "
$ws.Range("E4").Value = "  +0.29%  "
Set-TextValue "D5" "521.21"
$ws.Range("E5").Value = "  -1.21%  "
Set-TextValue "D6" "128.33"
$ws.Range("E6").Value = "  -5.10%  "
$ws.Range("E7").Value = "  -0.03%  "
Set-TextValue "D8" "0.560"
$ws.Range("E8").Value = "  -1.33%  "
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("E11").Value = "  -5.55%  "
Set-TextValue "D12" "0.320"
$ws.Range("E12").Value = "  -4.70%  "
$ws.Range("D13").Value = "2.876.35"
$ws.Range("D14").Value = "57.649.36"
$ws.Range("E14").Value = "  -2.38%  "
Set-TextValue "D15" "21.50"
$ws.Range("E15").Value = "  -3.90%  "
$ws.Range("E16").Value = "  -3.69%  "
$ws.Range("D17").Value = "2.443.34"
$ws.Range("E17").Value = "  -3.12%  "
Set-TextValue "D18" "10.36"
$ws.Range("E18").Value = "  -3.76%  "
$ws.Range("E19").Value = "  -2.79%  "
Set-TextValue "D20" "311.12"
$ws.Range("E20").Value = "  -4.00%  "
Set-TextValue "D21" "6.08"
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("E22").Value = "  +0.26%  "
Set-TextValue "D23" "64.82"
$ws.Range("E23").Value = "  -0.64%  "
Set-TextValue "D24" "0.399"
$ws.Range("E24").Value = "  -2.86%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "2.560.90"
$ws.Range("E25").Value = "  -3.60%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D26" "0.997"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -3.70%  "
$ws.Range("E28").Value = "  -4.22%  "
Set-TextValue "D29" "173.03"
$ws.Range("E29").Value = "  +1.95%  "
$ws.Range("E30").Value = "  -3.86%  "
$ws.Range("E31").Value = "  -3.32%  "
Set-TextValue "D32" "6.14"
$ws.Range("E32").Value = "  -3.51%  "
Set-TextValue "D33" "1.13"
$ws.Range("E33").Value = "  -8.73%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  -0.13%  "
Set-TextValue "D36" "17.78"
$ws.Range("E36").Value = "  -2.90%  "
$ws.Range("E37").Value = "  -7.40%  "
Set-TextValue "D38" "3.74"
$ws.Range("E38").Value = "  -6.38%  "
Set-TextValue "D39" "36.29"
$ws.Range("E39").Value = "  -1.18%  "
Set-TextValue "D40" "0.795"
$ws.Range("E40").Value = "  +1.05%  "
$ws.Range("E41").Value = "  -5.77%  "
$ws.Range("E42").Value = "  -3.70%  "
Set-TextValue "D43" "0.579"
$ws.Range("E43").Value = "  -3.89%  "
Set-TextValue "D44" "4.75"
$ws.Range("E44").Value = "  -7.09%  "
Set-TextValue "D45" "0.0917"
$ws.Range("E45").Value = "  -0.83%  "
Set-TextValue "D46" "252.54"
$ws.Range("E46").Value = "  -10.03%  "
Set-TextValue "D47" "119.71"
$ws.Range("E47").Value = "  -11.32%  "
Set-TextValue "D48" "0.0490"
$ws.Range("E48").Value = "  -3.63%  "
$ws.Range("E49").Value = "  -3.82%  "
$ws.Range("E50").Value = "  -5.39%  "
$ws.Range("E51").Value = "  -0.63%  "
